# Apply the "item points / skill points / buy cost" restructuring to
# the AttributeCardData sheet.
#
# Summary of the change:
#   - Insert a new column B named "rank" (maps to the new "skill points"
#     per-attribute value of 3) before the existing "maxCount" column
#     (now "item points" buy cost, stays 5), which pushes cardNameEn and
#     imageFile one column to the right.
#   - cardNameEn values change from lower-case png-ish words to proper
#     capitalised English attribute names.
#   - The "vigor.png" image filename becomes "Vitality.png" (and its
#     cardNameEn becomes "Vitality" instead of the old Chinese-derived
#     value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("maxCount"), shifting maxCount,
# cardNameEn, imageFile one column to the right.
$ws.Range("B1:B6").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "cardName"
$ws.Range("B1").Value = "rank"
$ws.Range("C1").Value = "maxCount"
$ws.Range("D1").Value = "cardNameEn"
$ws.Range("E1").Value = "imageFile"

# Row 2 - Strength / 力量
$ws.Range("A2").Value = "力量"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = "Strength"
$ws.Range("E2").Value = "strength.png"

# Row 3 - Vitality / 体质
$ws.Range("A3").Value = "体质"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "Vitality"
$ws.Range("E3").Value = "Vitality.png"

# Row 4 - Agility / 敏捷
$ws.Range("A4").Value = "敏捷"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "Agility"
$ws.Range("E4").Value = "agility.png"

# Row 5 - Awareness / 感知
$ws.Range("A5").Value = "感知"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "Awareness"
$ws.Range("E5").Value = "awareness.png"

# Row 6 - Intelligence / 智力
$ws.Range("A6").Value = "智力"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Intelligence"
$ws.Range("E6").Value = "intelligence.png"

# Selection moves to D7 in the saved file.
$ws.Range("D7").Select()

$wb.Save()
